# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2233
$ws1.Range("F4").Value = 13529
$ws1.Range("F5").Value = 82
$ws1.Range("F9").Value = 1200
$ws1.Range("F11").Value = 13853
$ws1.Range("F12").Value = 14619
$ws1.Range("F22").Value = 9
$ws1.Range("F26").Value = 5617
$ws1.Range("F29").Value = 5372
$ws1.Range("F31").Value = 36
$ws1.Range("F32").Value = 191

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2233
$ws4.Range("F4").Value = 13529
$ws4.Range("F5").Value = 82
$ws4.Range("F10").Value = 1200
$ws4.Range("F12").Value = 13853
$ws4.Range("F13").Value = 14619
$ws4.Range("F23").Value = 9
$ws4.Range("F27").Value = 5617
$ws4.Range("F30").Value = 5372
$ws4.Range("F32").Value = 36
$ws4.Range("F33").Value = 191
